$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize labels that become standalone headers (localization cleanup)
$ws.Range("C6").Value = "By sex"
$ws.Range("C12").Value = "By territory"

# Row 22 "Age" header block gets more descriptive text in 3 languages
$ws.Range("A22").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("B22").Value = "По возрасту (в годах)"
$ws.Range("C22").Value = "By age (in years)"

# Wealth quintile header capitalized
$ws.Range("C34").Value = "Wealth quintile"

# Row 27 + 30 (functional-difficulties sub-headers) get taller rows + wrap styling
$ws.Rows("27").RowHeight = 24
$ws.Range("A27:B27").WrapText = $true

$ws.Rows("30").RowHeight = 36
$ws.Range("A30:B30").WrapText = $true

# Placeholder "-" values for previously-blank data cells
$ws.Range("D28").Value = "-"
$ws.Range("D29").Value = "-"
$ws.Range("D31").Value = "-"
$ws.Range("D32").Value = "-"
$ws.Range("D33").Value = "-"

# Restore cursor/selection position
$ws.Range("B30").Select()
